$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the F column values that changed (row 1,2,3,5,9,10,11,12,13).
# Row 4 and 6,7,14,15 are formulas / unchanged and will recalc automatically.
$ws.Range("F1").Value = 0.25
$ws.Range("F2").Value = 0.7
$ws.Range("F3").Value = 0.9
$ws.Range("F5").Value = 0.92
$ws.Range("F9").Value = 0.92
$ws.Range("F10").Value = 0.94
$ws.Range("F11").Value = 0.96
$ws.Range("F12").Value = 0.97
$ws.Range("F13").Value = 0.98

# Update the selected cell/range on the sheet to match the new view.
$ws.Range("F2").Select()
